$d = $word.ActiveDocument

# --- Step 1: the document currently has a stray "_GoBack" bookmark sitting
#     right after "correta" (before the final period). That bookmark is being
#     relocated, so remove it from its old spot first.
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks.Item("_GoBack").Delete()
}

# --- Step 2: fix the typo "discursões" -> "discussões" in the heading
#     "Tópicos de discursões da 1ª Reunião:" (the "r" becomes an "s").
$heading = $d.Content
$heading.Find.Execute("Tópicos de discursões da 1ª Reunião:", $true, $false, $false, `
                       $false, $false, $true, 1, $false, "", 0)

if ($heading.Find.Found) {
    $headingStart = $heading.Start

    # Locate the single "r" in "discu[r]sões" relative to the match start and
    # turn it into an "s" (...discu| -> ...discus|sões...).
    $prefix = "Tópicos de discu"
    $rOffset = $prefix.Length
    $rChar = $d.Range($headingStart + $rOffset, $headingStart + $rOffset + 1)
    $rChar.Text = "s"

    # --- Step 3: drop the relocated "_GoBack" bookmark right after the fix,
    #     i.e. between "...discus" and "sões da 1ª Reunião:".
    $bmStart = $headingStart + $rOffset + 1
    $bmRange = $d.Range($bmStart, $bmStart)
    $d.Bookmarks.Add("_GoBack", $bmRange)
}
